$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "19-10="
$t.Cell(1,2).Range.Text = "5+70="
$t.Cell(1,3).Range.Text = "63-40="
$t.Cell(1,4).Range.Text = "40+13="
$t.Cell(1,5).Range.Text = "62-59="
$t.Cell(2,1).Range.Text = "76-29="
$t.Cell(2,2).Range.Text = "28+27="
$t.Cell(2,3).Range.Text = "64+30="
$t.Cell(2,4).Range.Text = "13+24="
$t.Cell(2,5).Range.Text = "26+29="
$t.Cell(3,1).Range.Text = "73-33="
$t.Cell(3,2).Range.Text = "72-34="
$t.Cell(3,3).Range.Text = "2+82="
$t.Cell(3,4).Range.Text = "46+32="
$t.Cell(3,5).Range.Text = "10+12="
$t.Cell(4,1).Range.Text = "63-38="
$t.Cell(4,2).Range.Text = "0+50="
$t.Cell(4,3).Range.Text = "21-7="
$t.Cell(4,4).Range.Text = "56+4="
$t.Cell(4,5).Range.Text = "40+22="
$t.Cell(5,1).Range.Text = "18-7="
$t.Cell(5,2).Range.Text = "94-82="
$t.Cell(5,3).Range.Text = "0+87="
$t.Cell(5,4).Range.Text = "67+19="
$t.Cell(5,5).Range.Text = "0+95="
$t.Cell(6,1).Range.Text = "26-14="
$t.Cell(6,2).Range.Text = "97-20="
$t.Cell(6,3).Range.Text = "27+32="
$t.Cell(6,4).Range.Text = "90-79="
$t.Cell(6,5).Range.Text = "28+53="
$t.Cell(7,1).Range.Text = "26+34="
$t.Cell(7,2).Range.Text = "74-50="
$t.Cell(7,3).Range.Text = "12+42="
$t.Cell(7,4).Range.Text = "53+28="
$t.Cell(7,5).Range.Text = "39+14="
$t.Cell(8,1).Range.Text = "84-11="
$t.Cell(8,2).Range.Text = "76-55="
$t.Cell(8,3).Range.Text = "24-22="
$t.Cell(8,4).Range.Text = "95-20="
$t.Cell(8,5).Range.Text = "77-54="
$t.Cell(9,1).Range.Text = "74-22="
$t.Cell(9,2).Range.Text = "40+19="
$t.Cell(9,3).Range.Text = "55-38="
$t.Cell(9,4).Range.Text = "39-1="
$t.Cell(9,5).Range.Text = "98-86="
$t.Cell(10,1).Range.Text = "85+14="
$t.Cell(10,2).Range.Text = "15+27="
$t.Cell(10,3).Range.Text = "17+51="
$t.Cell(10,4).Range.Text = "6+29="
$t.Cell(10,5).Range.Text = "29+48="
$t.Cell(11,1).Range.Text = "82-6="
$t.Cell(11,2).Range.Text = "37+45="
$t.Cell(11,3).Range.Text = "54-48="
$t.Cell(11,4).Range.Text = "70-28="
$t.Cell(11,5).Range.Text = "93-69="
$t.Cell(12,1).Range.Text = "83-58="
$t.Cell(12,2).Range.Text = "48+36="
$t.Cell(12,3).Range.Text = "50+39="
$t.Cell(12,4).Range.Text = "59-48="
$t.Cell(12,5).Range.Text = "55+35="
$t.Cell(13,1).Range.Text = "10+19="
$t.Cell(13,2).Range.Text = "28-7="
$t.Cell(13,3).Range.Text = "98-16="
$t.Cell(13,4).Range.Text = "61-35="
$t.Cell(13,5).Range.Text = "95-20="
$t.Cell(14,1).Range.Text = "24+70="
$t.Cell(14,2).Range.Text = "82+10="
$t.Cell(14,3).Range.Text = "26-19="
$t.Cell(14,4).Range.Text = "62+31="
$t.Cell(14,5).Range.Text = "19+4="
$t.Cell(15,1).Range.Text = "46+15="
$t.Cell(15,2).Range.Text = "53+32="
$t.Cell(15,3).Range.Text = "27+3="
$t.Cell(15,4).Range.Text = "81+13="
$t.Cell(15,5).Range.Text = "14+83="
$t.Cell(16,1).Range.Text = "43+22="
$t.Cell(16,2).Range.Text = "72+21="
$t.Cell(16,3).Range.Text = "60-45="
$t.Cell(16,4).Range.Text = "99-16="
$t.Cell(16,5).Range.Text = "87-25="
$t.Cell(17,1).Range.Text = "38-32="
$t.Cell(17,2).Range.Text = "67+29="
$t.Cell(17,3).Range.Text = "99-85="
$t.Cell(17,4).Range.Text = "47+30="
$t.Cell(17,5).Range.Text = "94-30="
$t.Cell(18,1).Range.Text = "59-7="
$t.Cell(18,2).Range.Text = "38+45="
$t.Cell(18,3).Range.Text = "63-39="
$t.Cell(18,4).Range.Text = "92-79="
$t.Cell(18,5).Range.Text = "30-10="
$t.Cell(19,1).Range.Text = "88-40="
$t.Cell(19,2).Range.Text = "63+24="
$t.Cell(19,3).Range.Text = "45+25="
$t.Cell(19,4).Range.Text = "66+8="
$t.Cell(19,5).Range.Text = "2+14="
$t.Cell(20,1).Range.Text = "83-12="
$t.Cell(20,2).Range.Text = "21+0="
$t.Cell(20,3).Range.Text = "40+43="
$t.Cell(20,4).Range.Text = "75-23="
$t.Cell(20,5).Range.Text = "20+20="
